$d = $word.ActiveDocument

# --- 1. Collapse the intro paragraph's five runs into a single run -----------
# The text content does not change, only the run structure (Word naturally
# coalesces runs with identical formatting when you do a Find/Replace over
# them). Scope the Find to just that paragraph's range so nothing else in
# the document can match.
$introPara = $d.Paragraphs(3)
$introRange = $introPara.Range
$introText = "I will be testing the 2 main user personas on 3 people. I will be giving each person a questionnaire to fill out in relations to the user persona so that I can get feedback. This is to ensure the user persona represents all standard users. "
$introRange.Find.Execute($introText, $true, $false, $false, $false, $false, $true, 1, $false, $introText, 2)

# --- 2. Remove the "Test User 2" and "Test User 3" heading paragraphs --------
$p2 = $null
$p3 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Test User 2`r") { $p2 = $p }
    if ($t -eq "Test User 3`r") { $p3 = $p }
}
if ($p2 -ne $null -and $p3 -ne $null) {
    $delRange = $d.Range($p2.Range.Start, $p3.Range.End)
    $delRange.Delete()
}
